$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1746987951807229
$ws.Range("C2").Value = 0.5993975903614458
$ws.Range("J2").Value = 0.003012048192771084
$ws.Range("P2").Value = 0.1506024096385542
$ws.Range("S2").Value = 0.07228915662650602
$ws.Range("B3").Value = 0.004807692307692308
$ws.Range("C3").Value = 0.02884615384615385
$ws.Range("J3").Value = 0.02884615384615385
$ws.Range("P3").Value = 0.8028846153846154
$ws.Range("S3").Value = 0.1346153846153846
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6808510638297872
$ws.Range("S4").Value = 0.2978723404255319
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("F6").Value = 0.05
$ws.Range("J6").Value = 0.3
$ws.Range("O6").Value = 0.008333333333333333
$ws.Range("Q6").Value = 0.1583333333333333
$ws.Range("R6").Value = 0.07916666666666666
$ws.Range("S6").Value = 0.3208333333333334
$ws.Range("B7").Value = 0.1263537906137184
$ws.Range("D7").Value = 0.01083032490974729
$ws.Range("F7").Value = 0.06137184115523465
$ws.Range("J7").Value = 0.1444043321299639
$ws.Range("O7").Value = 0.003610108303249098
$ws.Range("Q7").Value = 0.2021660649819494
$ws.Range("R7").Value = 0.09386281588447654
$ws.Range("S7").Value = 0.3574007220216607
$ws.Range("B8").Value = 0.1032388663967611
$ws.Range("D8").Value = 0.01417004048582996
$ws.Range("E8").Value = 0.004048582995951417
$ws.Range("F8").Value = 0.05060728744939271
$ws.Range("J8").Value = 0.1133603238866397
$ws.Range("O8").Value = 0.01619433198380567
$ws.Range("Q8").Value = 0.2024291497975708
$ws.Range("R8").Value = 0.07894736842105263
$ws.Range("S8").Value = 0.4170040485829959
$ws.Range("B9").Value = 0.09134615384615384
$ws.Range("D9").Value = 0.009615384615384616
$ws.Range("F9").Value = 0.01442307692307692
$ws.Range("J9").Value = 0.09134615384615384
$ws.Range("O9").Value = 0.01923076923076923
$ws.Range("Q9").Value = 0.2355769230769231
$ws.Range("R9").Value = 0.09134615384615384
$ws.Range("S9").Value = 0.4471153846153846
$ws.Range("B10").Value = 0.109118086696562
$ws.Range("D10").Value = 0.02316890881913303
$ws.Range("E10").Value = 0.001494768310911809
$ws.Range("F10").Value = 0.06801195814648729
$ws.Range("J10").Value = 0.1188340807174888
$ws.Range("O10").Value = 0.01420029895366218
$ws.Range("Q10").Value = 0.2085201793721973
$ws.Range("R10").Value = 0.08819133034379671
$ws.Range("S10").Value = 0.3684603886397608
$ws.Range("G11").Value = 0.1813842482100239
$ws.Range("J11").Value = 0.08353221957040573
$ws.Range("K11").Value = 0.2147971360381861
$ws.Range("L11").Value = 0.5083532219570406
$ws.Range("S11").Value = 0.01193317422434368
$ws.Range("G12").Value = 0.7207207207207207
$ws.Range("J12").Value = 0.1801801801801802
$ws.Range("K12").Value = 0.02252252252252252
$ws.Range("L12").Value = 0.03603603603603604
$ws.Range("S12").Value = 0.04054054054054054
$ws.Range("G13").Value = 0.7903225806451613
$ws.Range("J13").Value = 0.1935483870967742
$ws.Range("S13").Value = 0.01612903225806452
$ws.Range("F15").Value = 0.0186046511627907
$ws.Range("H15").Value = 0.1627906976744186
$ws.Range("I15").Value = 0.06976744186046512
$ws.Range("J15").Value = 0.3627906976744186
$ws.Range("K15").Value = 0.08837209302325581
$ws.Range("M15").Value = 0.01395348837209302
$ws.Range("O15").Value = 0.06976744186046512
$ws.Range("S15").Value = 0.213953488372093
$ws.Range("F16").Value = 0.0125
$ws.Range("H16").Value = 0.1208333333333333
$ws.Range("I16").Value = 0.08749999999999999
$ws.Range("J16").Value = 0.4041666666666667
$ws.Range("K16").Value = 0.1458333333333333
$ws.Range("M16").Value = 0.02916666666666667
$ws.Range("O16").Value = 0.05
$ws.Range("S16").Value = 0.15
$ws.Range("F17").Value = 0.02702702702702703
$ws.Range("H17").Value = 0.1718146718146718
$ws.Range("I17").Value = 0.09266409266409266
$ws.Range("J17").Value = 0.3783783783783784
$ws.Range("K17").Value = 0.1332046332046332
$ws.Range("M17").Value = 0.02895752895752896
$ws.Range("N17").Value = 0.003861003861003861
$ws.Range("O17").Value = 0.05791505791505792
$ws.Range("S17").Value = 0.1061776061776062
$ws.Range("F18").Value = 0.02727272727272727
$ws.Range("H18").Value = 0.1545454545454545
$ws.Range("I18").Value = 0.08636363636363636
$ws.Range("J18").Value = 0.3863636363636364
$ws.Range("K18").Value = 0.1454545454545454
$ws.Range("M18").Value = 0.03181818181818181
$ws.Range("O18").Value = 0.08181818181818182
$ws.Range("S18").Value = 0.08636363636363636
$ws.Range("F19").Value = 0.02462686567164179
$ws.Range("H19").Value = 0.2291044776119403
$ws.Range("I19").Value = 0.0791044776119403
$ws.Range("J19").Value = 0.3440298507462687
$ws.Range("K19").Value = 0.1253731343283582
$ws.Range("M19").Value = 0.02313432835820895
$ws.Range("N19").Value = 0.0007462686567164179
$ws.Range("O19").Value = 0.06044776119402985
$ws.Range("S19").Value = 0.1134328358208955
